$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string (shared string used by A1)
$ws.Range("A1").Value = "Datos actualizados a 6 de Abril de 2020 a las 14:22"

# Row 13
$ws.Range("B13").Value = 21652
$ws.Range("C13").Value = 552
$ws.Range("E13").Value = 13620

# Row 17
$ws.Range("B17").Value = 12162
$ws.Range("C17").Value = 111
$ws.Range("E17").Value = 8479

# Row 19
$ws.Range("B19").Value = 11298
$ws.Range("C19").Value = 44
$ws.Range("E19").Value = 10682
$ws.Range("G19").Value = 3
$ws.Range("H19").Value = 489

# Row 22
$ws.Range("B22").Value = 7206
$ws.Range("C22").Value = 376
$ws.Range("E22").Value = 6524
$ws.Range("F22").Value = 590
$ws.Range("G22").Value = 76
$ws.Range("H22").Value = 477

# Row 24
$ws.Range("E24").Value = 3322
$ws.Range("G24").Value = 4
$ws.Range("H24").Value = 41

# Row 52
$ws.Range("E52").Value = 1226
$ws.Range("G52").Value = 2
$ws.Range("H52").Value = 48

# Row 73
$ws.Range("E73").Value = 592
$ws.Range("G73").Value = 5
$ws.Range("H73").Value = 28

# Row 76
$ws.Range("F76").Value = 11

# Row 113
$ws.Range("B113").Value = 178
$ws.Range("C113").Value = 2
$ws.Range("E113").Value = 139

# Row 115
$ws.Range("A115").Value = "Kenia"
$ws.Range("B115").Value = 158
$ws.Range("C115").Value = 16
$ws.Range("D115").Value = 4
$ws.Range("E115").Value = 150
$ws.Range("F115").Value = 2
$ws.Range("H115").Value = 4

# Row 116
$ws.Range("A116").Value = "Consejo Danes para los Refugiados"
$ws.Range("B116").Value = 154
$ws.Range("D116").Value = 3
$ws.Range("E116").Value = 133
$ws.Range("F116").Value = 0
$ws.Range("H116").Value = 18

# Row 117
$ws.Range("A117").Value = "Martinica"
$ws.Range("B117").Value = 149
$ws.Range("D117").Value = 50
$ws.Range("E117").Value = 95
$ws.Range("F117").Value = 21
$ws.Range("H117").Value = 4

# Row 118
$ws.Range("A118").Value = "Mayotte"
$ws.Range("B118").Value = 147
$ws.Range("D118").Value = 14
$ws.Range("E118").Value = 131
$ws.Range("F118").Value = 3
$ws.Range("H118").Value = 2

